# Scheduled runner update: refresh market-price/profit figures pulled from
# the external market API for several leve rows across the item-crafting
# sheets (ALC, ARM, BSM, CRP, CUL, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 146.1
$ws.Range("I9").Value = 157.625
$ws.Range("J9").Value = 100
$ws.Range("K9").Value = 157.625
$ws.Range("L9").Value = 100
$ws.Range("M9").Value = 11.375
$ws.Range("N9").Value = -438

$ws.Range("H18").Value = 593.8570999999999
$ws.Range("I18").Value = 593.8570999999999
$ws.Range("K18").Value = 593.8570999999999
$ws.Range("M18").Value = -309.8570999999999

$ws.Range("H62").Value = 2138.3333
$ws.Range("I62").Value = 2066
$ws.Range("J62").Value = 2500
$ws.Range("K62").Value = 2066
$ws.Range("L62").Value = 2500
$ws.Range("M62").Value = -1442
$ws.Range("N62").Value = -3748

$ws.Range("H65").Value = 2138.3333
$ws.Range("I65").Value = 2066
$ws.Range("J65").Value = 2500
$ws.Range("K65").Value = 10330
$ws.Range("L65").Value = 12500
$ws.Range("M65").Value = -7210
$ws.Range("N65").Value = -18740

$ws.Range("H100").Value = 2186.7273
$ws.Range("I100").Value = 2264
$ws.Range("J100").Value = 2094
$ws.Range("K100").Value = 2264
$ws.Range("L100").Value = 2094
$ws.Range("M100").Value = -1723
$ws.Range("N100").Value = -3176

$ws.Range("H103").Value = 796.53845
$ws.Range("I103").Value = 648.75
$ws.Range("K103").Value = 1946.25
$ws.Range("M103").Value = -1360.25

$ws.Range("H132").Value = 5819096.5
$ws.Range("I132").Value = 5819096.5
$ws.Range("K132").Value = 17457289.5
$ws.Range("M132").Value = -17454759.5

$ws.Range("H137").Value = 1968.5927
$ws.Range("I137").Value = 1242
$ws.Range("K137").Value = 3726
$ws.Range("M137").Value = -1176

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 38539.777
$ws.Range("I2").Value = 1061.2222
$ws.Range("J2").Value = 113496.89
$ws.Range("K2").Value = 1061.2222
$ws.Range("L2").Value = 113496.89
$ws.Range("M2").Value = -948.2221999999999
$ws.Range("N2").Value = -113722.89

$ws.Range("H61").Value = 1501.7273
$ws.Range("I61").Value = 761.2646999999999
$ws.Range("J61").Value = 2700.5715
$ws.Range("K61").Value = 761.2646999999999
$ws.Range("L61").Value = 2700.5715
$ws.Range("M61").Value = -549.2646999999999
$ws.Range("N61").Value = -3124.5715

$ws.Range("H116").Value = 38539.777
$ws.Range("I116").Value = 1061.2222
$ws.Range("J116").Value = 113496.89
$ws.Range("K116").Value = 1061.2222
$ws.Range("L116").Value = 113496.89
$ws.Range("M116").Value = 1232.7778
$ws.Range("N116").Value = -118084.89

$ws.Range("H136").Value = 1501.7273
$ws.Range("I136").Value = 761.2646999999999
$ws.Range("J136").Value = 2700.5715
$ws.Range("K136").Value = 2283.7941
$ws.Range("L136").Value = 8101.7145
$ws.Range("M136").Value = 266.2058999999999
$ws.Range("N136").Value = -13201.7145

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 38539.777
$ws.Range("I3").Value = 1061.2222
$ws.Range("J3").Value = 113496.89
$ws.Range("K3").Value = 1061.2222
$ws.Range("L3").Value = 113496.89
$ws.Range("M3").Value = -947.2221999999999
$ws.Range("N3").Value = -113724.89

$ws.Range("H99").Value = 1329.05
$ws.Range("I99").Value = 1189.091
$ws.Range("J99").Value = 1500.1111
$ws.Range("K99").Value = 1189.091
$ws.Range("L99").Value = 1500.1111
$ws.Range("M99").Value = 308.9090000000001
$ws.Range("N99").Value = -4496.1111

$ws.Range("H107").Value = 16689649
$ws.Range("I107").Value = 20860804
$ws.Range("K107").Value = 20860804
$ws.Range("M107").Value = -20858884

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H32").Value = 7000
$ws.Range("J32").Value = 7000
$ws.Range("L32").Value = 7000
$ws.Range("N32").Value = -7632

$ws.Range("H107").Value = 4351.6665
$ws.Range("I107").Value = 7805.4287
$ws.Range("J107").Value = 632.2308
$ws.Range("K107").Value = 7805.4287
$ws.Range("L107").Value = 632.2308
$ws.Range("M107").Value = -5885.4287
$ws.Range("N107").Value = -4472.2308

$ws.Range("H122").Value = 2259
$ws.Range("I122").Value = 2192
$ws.Range("J122").Value = 2476.75
$ws.Range("K122").Value = 6576
$ws.Range("L122").Value = 7430.25
$ws.Range("M122").Value = -4126
$ws.Range("N122").Value = -12330.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 543.25
$ws.Range("I14").Value = 543.25
$ws.Range("K14").Value = 1629.75
$ws.Range("M14").Value = -1456.75

$ws.Range("H17").Value = 2091.923
$ws.Range("I17").Value = 2250
$ws.Range("J17").Value = 1956.4286
$ws.Range("K17").Value = 6750
$ws.Range("L17").Value = 5869.2858
$ws.Range("M17").Value = -6581
$ws.Range("N17").Value = -6207.2858

$ws.Range("H34").Value = 1581.7273
$ws.Range("J34").Value = 1899.8889
$ws.Range("L34").Value = 5699.6667
$ws.Range("N34").Value = -5867.6667

$ws.Range("H39").Value = 31800
$ws.Range("J39").Value = 31800
$ws.Range("L39").Value = 95400
$ws.Range("N39").Value = -95988

$ws.Range("H44").Value = 800
$ws.Range("I44").Value = 600
$ws.Range("K44").Value = 1800
$ws.Range("M44").Value = -1402

$ws.Range("H55").Value = 9208.733
$ws.Range("J55").Value = 9823.643
$ws.Range("L55").Value = 29470.929
$ws.Range("N55").Value = -29824.929

$ws.Range("H107").Value = 262987.84
$ws.Range("I107").Value = 365.66666
$ws.Range("J107").Value = 819128.9399999999
$ws.Range("K107").Value = 1096.99998
$ws.Range("L107").Value = 2457386.82
$ws.Range("M107").Value = 823.0000199999999
$ws.Range("N107").Value = -2461226.82

$ws.Range("H113").Value = 1190.6666
$ws.Range("I113").Value = 528.3077
$ws.Range("J113").Value = 1697.1765
$ws.Range("K113").Value = 1584.9231
$ws.Range("L113").Value = 5091.529500000001
$ws.Range("M113").Value = 585.0769
$ws.Range("N113").Value = -9431.529500000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").ClearContents()

$ws.Range("H20").Value = 26752.25
$ws.Range("J20").Value = 26752.25
$ws.Range("L20").Value = 26752.25
$ws.Range("N20").Value = -27204.25

$ws.Range("H22").Value = 471.64285
$ws.Range("I22").Value = 367
$ws.Range("J22").Value = 500.18182
$ws.Range("K22").Value = 367
$ws.Range("L22").Value = 500.18182
$ws.Range("M22").Value = -72
$ws.Range("N22").Value = -1090.18182

$ws.Range("H27").Value = 471.64285
$ws.Range("I27").Value = 367
$ws.Range("J27").Value = 500.18182
$ws.Range("K27").Value = 367
$ws.Range("L27").Value = 500.18182
$ws.Range("M27").Value = -260
$ws.Range("N27").Value = -714.18182

$ws.Range("H61").Value = 2930
$ws.Range("I61").Value = 2500
$ws.Range("K61").Value = 2500
$ws.Range("M61").Value = -2298

$ws.Range("H113").Value = 2930
$ws.Range("I113").Value = 2500
$ws.Range("K113").Value = 2500
$ws.Range("M113").Value = -330

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 6793.6665
$ws.Range("J49").Value = 6793.6665
$ws.Range("L49").Value = 6793.6665
$ws.Range("N49").Value = -7253.6665

$ws.Range("H54").Value = 6944.25
$ws.Range("J54").Value = 6944.25
$ws.Range("L54").Value = 6944.25
$ws.Range("N54").Value = -7984.25

$ws.Range("H81").Value = 500974.5
$ws.Range("I81").Value = 334633
$ws.Range("K81").Value = 669266
$ws.Range("M81").Value = -668205

$ws.Range("H84").Value = 500974.5
$ws.Range("I84").Value = 334633
$ws.Range("K84").Value = 3346330
$ws.Range("M84").Value = -3341026

$ws.Range("H107").Value = 59835.105
$ws.Range("I107").Value = 13047.444
$ws.Range("J107").Value = 101944
$ws.Range("K107").Value = 39142.33199999999
$ws.Range("L107").Value = 305832
$ws.Range("M107").Value = -37222.33199999999
$ws.Range("N107").Value = -309672
